# CpiChartRepository.py - rewrite the month-header row (row 1, columns B..AW)
# from "2021年1月" style labels to compact "202101" style labels, so the new
# chart-drawing code can parse/sort them as YYYYMM strings.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet 1")

$startYear = 2021
$startMonth = 1
$startCol = 2   # column B is the first month column; A holds the row label

for ($i = 0; $i -lt 48; $i++) {
    $year = $startYear + [math]::Floor(($startMonth - 1 + $i) / 12)
    $month = (($startMonth - 1 + $i) % 12) + 1
    $label = "{0}{1:D2}" -f $year, $month

    $cell = $ws.Cells.Item(1, $startCol + $i)
    # Leading apostrophe forces the numeric-looking text to be stored as a
    # literal string (e.g. "202101"), not coerced into the number 202101.
    $cell.Value = "'" + $label
}
